$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the value "Test" into A1 (becomes a shared string).
$ws.Range("A1").Value = "Test"

# Leave the selection on B3, matching the saved view state.
$ws.Range("B3").Select() | Out-Null

# Switch the sheet's print orientation to Portrait (xlPortrait = 1).
$ws.PageSetup.Orientation = 1
